$wb = $excel.ActiveWorkbook

# Sheets
$wsIndex = $wb.Worksheets.Item("Model Index")
$wsIngreso = $wb.Worksheets.Item("Facts - Ingreso")
$wsLocal = $wb.Worksheets.Item("Facts - Localizacion")

# Data change: mark "Localizacion" fact group as also related to PRODUCTO
$wsIndex.Range("C5").Value = "X"

# Selections per sheet (view state)
$wsIndex.Range("D10").Select()
$wsIngreso.Range("B4").Select()
$wsLocal.Range("B3").Select()

# Active sheet should be "Facts - Localizacion" (was "Model Index")
$wsLocal.Activate()
